$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 2 and row 3 and need to be swapped:
# A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $val2 = $cellRow2.Value2
    $val3 = $cellRow3.Value2

    $cellRow2.Value2 = $val3
    $cellRow3.Value2 = $val2
}
